$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Split the "Summarize<br/>What Finance Teams Achieve..." paragraph into
#    two separate paragraphs: a standalone "Summarize" Heading3 paragraph
#    (styled with the Heading1Char character style) followed by the original
#    callout-box heading (minus the "Summarize" run).
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Summarize", $true)
$para = $rng.Paragraphs(1).Range

$xml1 = '<w:p ' + $ns + '>' +
          '<w:pPr>' +
            '<w:pStyle w:val="Heading3"/>' +
            '<w:spacing w:before="281" w:after="281"/>' +
          '</w:pPr>' +
          '<w:r><w:t>Summarize</w:t></w:r>' +
        '</w:p>' +
        '<w:p ' + $ns + '>' +
          '<w:pPr>' +
            '<w:pStyle w:val="Heading3"/>' +
            '<w:spacing w:before="281" w:after="281"/>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:highlight w:val="yellow"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
          '<w:r><w:br/></w:r>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
            '</w:rPr>' +
            '<w:t>What Finance Teams Achieve with Synoptix AI (Callout box)</w:t>' +
          '</w:r>' +
        '</w:p>'
$para.InsertXML($xml1)

# Set the paragraph mark's run properties (w:pPr/w:rPr) to reference the
# Heading1Char character style - InsertXML preserves w:rStyle when it is
# inside w:pPr/w:rPr, just not inside a run's own w:rPr.
$rng = $d.Content
[void]$rng.Find.Execute("Summarize", $true)
$summarizePara = $rng.Paragraphs(1).Range
$pPrXml = '<w:p ' + $ns + '>' +
            '<w:pPr>' +
              '<w:pStyle w:val="Heading3"/>' +
              '<w:spacing w:before="281" w:after="281"/>' +
              '<w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr>' +
            '</w:pPr>' +
            '<w:r><w:t>Summarize</w:t></w:r>' +
          '</w:p>'
$summarizePara.InsertXML($pPrXml)

# Re-apply the Heading1Char character style to the "Summarize" run itself -
# InsertXML silently drops w:rStyle on a run's own w:rPr, so we set it
# through the Range.Style property instead (by display name). This must be
# the last edit made to this particular run/paragraph.
$rng = $d.Content
[void]$rng.Find.Execute("Summarize", $true)
$summarizePara = $rng.Paragraphs(1).Range
$runRng = $d.Range($summarizePara.Start, $summarizePara.Start + 9)
$runRng.Style = "Heading 1 Char"

# ---------------------------------------------------------------------------
# 2) Merge the two runs of the "Finance teams drive strategic..." paragraph
#    back into a single run (drops the mid-sentence lastRenderedPageBreak).
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("them down. Synoptix", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "them down. Synoptix", 2)

# ---------------------------------------------------------------------------
# 3) Split the "Month-end shouldn't feel like a scramble..." paragraph into
#    two runs, inserting a lastRenderedPageBreak right before "variance
#    reports".
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Month-end", $true)
$monthEndPara = $rng.Paragraphs(1).Range

$xml2 = '<w:p ' + $ns + '>' +
          '<w:pPr><w:spacing w:before="240" w:after="240"/></w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">Month-end shouldn' + [char]0x2019 + 't feel like a scramble. Synoptix supports automation across recurring financial tasks: pulling ledger entries, aggregating expense data, generating </w:t>' +
          '</w:r>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
            '</w:rPr>' +
            '<w:lastRenderedPageBreak/>' +
            '<w:t>variance reports, and preparing board-ready summaries. Less time on manual work means more focus on forecasting, investment analysis, and strategic planning.</w:t>' +
          '</w:r>' +
        '</w:p>'
$monthEndPara.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 4) Remove the lastRenderedPageBreak before "Support " (now precedes the
#    "Month-end" split instead).
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Support the Business with Fast, Accurate Financial Guidance", $true)
$supportPara = $rng.Paragraphs(1).Range

$xml3 = '<w:p ' + $ns + '>' +
          '<w:pPr><w:pStyle w:val="Heading2"/><w:spacing w:before="360"/></w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:sz w:val="34"/><w:szCs w:val="34"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">Support </w:t>' +
          '</w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:sz w:val="34"/><w:szCs w:val="34"/>' +
            '</w:rPr>' +
            '<w:t>the Business</w:t>' +
          '</w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:sz w:val="34"/><w:szCs w:val="34"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve"> with Fast, Accurate Financial Guidance</w:t>' +
          '</w:r>' +
        '</w:p>'
$supportPara.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 5) Add a lastRenderedPageBreak before "Advanced Security and Governance".
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Advanced Security and Governance", $true)
$advSecPara = $rng.Paragraphs(1).Range

$xml4 = '<w:p ' + $ns + '>' +
          '<w:pPr><w:pStyle w:val="Heading3"/><w:spacing w:before="280"/></w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:sz w:val="26"/><w:szCs w:val="26"/>' +
            '</w:rPr>' +
            '<w:lastRenderedPageBreak/>' +
            '<w:t>Advanced Security and Governance</w:t>' +
          '</w:r>' +
        '</w:p>'
$advSecPara.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 6) Remove the lastRenderedPageBreak before "Finance Enablement That
#    Scales".
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("Finance Enablement That Scales", $true)
$enablementPara = $rng.Paragraphs(1).Range

$xml5 = '<w:p ' + $ns + '>' +
          '<w:pPr><w:pStyle w:val="Heading3"/><w:spacing w:before="280"/></w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '<w:b/><w:bCs/>' +
              '<w:color w:val="000000" w:themeColor="text1"/>' +
              '<w:sz w:val="26"/><w:szCs w:val="26"/>' +
            '</w:rPr>' +
            '<w:t>Finance Enablement That Scales</w:t>' +
          '</w:r>' +
        '</w:p>'
$enablementPara.InsertXML($xml5)

Write-Output "done"
